$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("kristina+literature_info")

# Update the latitude for bcitwo (row 3) to 45
$ws.Range("L3").Value = 45

# Update the lat_zone for bcitwo (row 3) to "temperate"
$ws.Range("N3").Value = "temperate"

# Move selection to N3, matching the saved workbook state
$ws.Range("N3").Select()
